$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 10:50"

# Row 17: Austria
$ws.Range("A17").Value = "Austria"
$ws.Range("B17").Value = 11224
$ws.Range("C17").Value = 95
$ws.Range("D17").Value = 2022
$ws.Range("E17").Value = 9044
$ws.Range("F17").Value = 227
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 158

# Row 33: Filipinas
$ws.Range("A33").Value = "Filipinas"
$ws.Range("B33").Value = 3018
$ws.Range("C33").Value = 385
$ws.Range("D33").Value = 52
$ws.Range("E33").Value = 2830
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 29
$ws.Range("H33").Value = 136

# Row 34: Rumania
$ws.Range("A34").Value = "Rumania"
$ws.Range("B34").Value = 2738
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 267
$ws.Range("E34").Value = 2355
$ws.Range("F34").Value = 78
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 116

# Row 39: Indonesia
$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 1986
$ws.Range("C39").Value = 196
$ws.Range("D39").Value = 134
$ws.Range("E39").Value = 1671
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 11
$ws.Range("H39").Value = 181

# Row 40: Tailandia
$ws.Range("A40").Value = "Tailandia"
$ws.Range("B40").Value = 1978
$ws.Range("C40").Value = 103
$ws.Range("D40").Value = 581
$ws.Range("E40").Value = 1378
$ws.Range("F40").Value = 23
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 19

# Row 41: Arabia Saudita
$ws.Range("A41").Value = "Arabia Saudita"
$ws.Range("B41").Value = 1885
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 328
$ws.Range("E41").Value = 1536
$ws.Range("F41").Value = 31
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 21

# Row 57: Estonia
$ws.Range("A57").Value = "Estonia"
$ws.Range("B57").Value = 961
$ws.Range("C57").Value = 103
$ws.Range("D57").Value = 48
$ws.Range("E57").Value = 901
$ws.Range("F57").Value = 16
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 12

# Row 58: Catar
$ws.Range("A58").Value = "Catar"
$ws.Range("B58").Value = 949
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 72
$ws.Range("E58").Value = 874
$ws.Range("F58").Value = 37
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 3

# Row 59: Ucrania
$ws.Range("A59").Value = "Ucrania"
$ws.Range("B59").Value = 942
$ws.Range("C59").Value = 45
$ws.Range("D59").Value = 19
$ws.Range("E59").Value = 900
$ws.Range("F59").Value = 16
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 23

# Row 60: Eslovenia
$ws.Range("A60").Value = "Eslovenia"
$ws.Range("B60").Value = 897
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 70
$ws.Range("E60").Value = 810
$ws.Range("F60").Value = 31
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 17

# Row 61: Nueva Zelanda
$ws.Range("A61").Value = "Nueva Zelanda"
$ws.Range("B61").Value = 868
$ws.Range("C61").Value = 71
$ws.Range("D61").Value = 103
$ws.Range("E61").Value = 764
$ws.Range("F61").Value = 2
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 1

# Row 62: Egipto
$ws.Range("A62").Value = "Egipto"
$ws.Range("B62").Value = 865
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 201
$ws.Range("E62").Value = 606
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 58

# Row 78: Azerbaiyan
$ws.Range("A78").Value = "Azerbaiyan"
$ws.Range("B78").Value = 443
$ws.Range("C78").Value = 43
$ws.Range("D78").Value = 32
$ws.Range("E78").Value = 406
$ws.Range("F78").Value = 7
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 5

# Row 79: Principado de Andorra
$ws.Range("A79").Value = "Principado de Andorra"
$ws.Range("B79").Value = 428
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 10
$ws.Range("E79").Value = 403
$ws.Range("F79").Value = 12
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 15

# Row 80: Eslovaquia
$ws.Range("A80").Value = "Eslovaquia"
$ws.Range("B80").Value = 426
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 5
$ws.Range("E80").Value = 420
$ws.Range("F80").Value = 3
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 1

# Row 81: Kuwait
$ws.Range("A81").Value = "Kuwait"
$ws.Range("B81").Value = 417
$ws.Range("C81").Value = 75
$ws.Range("D81").Value = 82
$ws.Range("E81").Value = 335
$ws.Range("F81").Value = 16
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 0

# Row 82: Costa Rica
$ws.Range("A82").Value = "Costa Rica"
$ws.Range("B82").Value = 396
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 6
$ws.Range("E82").Value = 388
$ws.Range("F82").Value = 11
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 2

# Row 83: Republica de Macedonia
$ws.Range("A83").Value = "Republica de Macedonia"
$ws.Range("B83").Value = 384
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 17
$ws.Range("E83").Value = 356
$ws.Range("F83").Value = 8
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 11

# Row 84: Uruguay
$ws.Range("A84").Value = "Uruguay"
$ws.Range("B84").Value = 369
$ws.Range("C84").Value = 19
$ws.Range("D84").Value = 68
$ws.Range("E84").Value = 297
$ws.Range("F84").Value = 13
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 4

# Row 85: Republica de Chipre
$ws.Range("A85").Value = "Republica de Chipre"
$ws.Range("B85").Value = 356
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 28
$ws.Range("E85").Value = 318
$ws.Range("F85").Value = 11
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 10

# Row 86: Taiwan
$ws.Range("A86").Value = "Taiwan"
$ws.Range("B86").Value = 348
$ws.Range("C86").Value = 9
$ws.Range("D86").Value = 50
$ws.Range("E86").Value = 293
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 5

# Row 106: Estado de Palestina
$ws.Range("A106").Value = "Estado de Palestina"
$ws.Range("B106").Value = 171
$ws.Range("C106").Value = 10
$ws.Range("D106").Value = 18
$ws.Range("E106").Value = 152
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 1

# Row 107: Mauricio
$ws.Range("A107").Value = "Mauricio"
$ws.Range("B107").Value = 169
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 0
$ws.Range("E107").Value = 162
$ws.Range("F107").Value = 1
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 7

# Row 177: Groenlandia
$ws.Range("A177").Value = "Groenlandia"
$ws.Range("B177").Value = 10
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 3
$ws.Range("E177").Value = 7
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 0

# Row 184: Republica de Africa Central
$ws.Range("A184").Value = "Republica de Africa Central"
$ws.Range("B184").Value = 8
$ws.Range("C184").Value = 5
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 8
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

# Row 185: Republica del Chad
$ws.Range("A185").Value = "Republica del Chad"
$ws.Range("B185").Value = 8
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 8
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 0

# Row 188: Fiyi
$ws.Range("A188").Value = "Fiyi"
$ws.Range("B188").Value = 7
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 0
$ws.Range("E188").Value = 7
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

# Row 189: Santa Sede
$ws.Range("A189").Value = "Santa Sede"
$ws.Range("B189").Value = 7
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 7
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 0

# Row 190: Nepal
$ws.Range("A190").Value = "Nepal"
$ws.Range("B190").Value = 7
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 1
$ws.Range("E190").Value = 6
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

# Row 191: Liberia
$ws.Range("A191").Value = "Liberia"
$ws.Range("B191").Value = 6
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 6
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 0

# Row 192: Cabo Verde
$ws.Range("A192").Value = "Cabo Verde"
$ws.Range("B192").Value = 6
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 0
$ws.Range("E192").Value = 5
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 1
